$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 623, shifting existing rows 623:664 down to 624:665
$ws.Rows.Item(623).Insert()

# Populate the newly inserted row 623 with the new data.
# Column A holds a date-formatted-as-text string ("2026/01/14"); force it to
# stay plain text (avoid Excel's automatic date-serial coercion) by
# temporarily marking the cell as Text, then stripping the format back off
# so the cell ends up with the default style (matching the rest of the sheet).
$ws.Cells.Item(623, 1).NumberFormat = "@"
$ws.Cells.Item(623, 1).Value = "2026/01/14"
$ws.Cells.Item(623, 1).ClearFormats()

$ws.Cells.Item(623, 2).Value = "水"
$ws.Cells.Item(623, 3).Value = 17
$ws.Cells.Item(623, 4).Value = 27
